# lophoc.xlsx edit: "decision tree" (cây quyết định) column C outcomes
# were renamed from có/không (yes/no) to nghỉ/học (rest/study), and the
# active sheet + selection moved from Sheet1 back to the "data" sheet.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("data")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# Column C on "data": replace "có" -> "nghỉ" and "không" -> "học"
$wsData.Range("C2").Value = "nghỉ"
$wsData.Range("C3").Value = "nghỉ"
$wsData.Range("C4").Value = "nghỉ"
$wsData.Range("C5").Value = "học"
$wsData.Range("C6").Value = "học"
$wsData.Range("C7").Value = "học"
$wsData.Range("C8").Value = "nghỉ"

# "Sheet1" keeps its own prior selection, just no longer the active tab
$wsSheet1.Range("C8").Select()

# Make "data" the active/selected sheet again, with D3 selected
$wsData.Activate()
$wsData.Range("D3").Select()
